$d = $word.ActiveDocument

$replacements = @(
    @("37×21=", "57×14="),
    @("41×27=", "34×37="),
    @("68×93=", "26×33="),
    @("29×13=", "51×88="),
    @("62×99=", "63×71="),
    @("23×51=", "11×68="),
    @("66×85=", "31×21="),
    @("41×29=", "73×24="),
    @("33×65=", "33×66="),
    @("52×17=", "92×74="),
    @("47×36=", "35×57="),
    @("68×51=", "43×37="),
    @("61×33=", "79×64="),
    @("24×99=", "87×46="),
    @("51×78=", "20×90="),
    @("40×64=", "35×96="),
    @("44×62=", "73×13="),
    @("17×42=", "92×50="),
    @("43×99=", "38×29="),
    @("64×47=", "19×76="),
    @("18×68=", "49×71="),
    @("69×32=", "20×45="),
    @("41×20=", "27×12="),
    @("18×35=", "49×91="),
    @("76×35=", "84×27=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
